# fix: consistent sort on championship log
#
# The "Log Campeonato" sheet lists, for each category, every player who
# scored championship points, grouped in blocks (one block per category,
# separated by a blank row). Each block must be sorted by championship
# points (column D) descending, and ties broken alphabetically by the
# player's name (column A) ascending. The data had drifted out of that
# order; this re-sorts every block in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log Campeonato Torneo 02 2022")

$startRow = 5
$maxScan = 1000

# --- 1. discover the category blocks (contiguous non-blank row runs) ---
$row = $startRow
$blankStreak = 0
$blockStart = -1
$blocks = @()

while ($row -le $maxScan) {
    $aVal = $ws.Cells.Item($row, 1).Value2
    $isBlank = ($aVal -eq $null) -or ($aVal -eq "")
    if (-not $isBlank) {
        if ($blockStart -eq -1) { $blockStart = $row }
        $blankStreak = 0
    } else {
        if ($blockStart -ne -1) {
            $blocks += , @($blockStart, $row - 1)
            $blockStart = -1
        }
        $blankStreak = $blankStreak + 1
        if ($blankStreak -ge 2) {
            break
        }
    }
    $row = $row + 1
}
if ($blockStart -ne -1) {
    $blocks += , @($blockStart, $row - 1)
}

# --- 2. sort each block in place: Points (D) desc, Player (A) asc ---
foreach ($block in $blocks) {
    $first = $block[0]
    $last = $block[1]
    $n = $last - $first + 1

    $items = @()
    for ($r = $first; $r -le $last; $r++) {
        $items += [PSCustomObject]@{
            A = $ws.Cells.Item($r, 1).Value2
            B = $ws.Cells.Item($r, 2).Value2
            C = $ws.Cells.Item($r, 3).Value2
            D = $ws.Cells.Item($r, 4).Value2
        }
    }

    # stable insertion sort: Points desc, then Player name asc
    for ($i = 1; $i -lt $n; $i++) {
        $key = $items[$i]
        $j = $i - 1
        while ($j -ge 0 -and (
            ($items[$j].D -lt $key.D) -or
            ($items[$j].D -eq $key.D -and $items[$j].A -gt $key.A)
        )) {
            $items[$j + 1] = $items[$j]
            $j = $j - 1
        }
        $items[$j + 1] = $key
    }

    for ($i = 0; $i -lt $n; $i++) {
        $r = $first + $i
        $ws.Cells.Item($r, 1).Value2 = $items[$i].A
        $ws.Cells.Item($r, 2).Value2 = $items[$i].B
        $ws.Cells.Item($r, 3).Value2 = $items[$i].C
        $ws.Cells.Item($r, 4).Value2 = $items[$i].D
    }
}
